$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textCells = @("D14", "D16", "D18", "D19", "D23", "D25", "D26", "D27", "D34", "D35", "D37", "D42", "D43", "D45", "D46", "D48")
foreach ($cellref in $textCells) {
    $ws.Range($cellref).NumberFormat = "@"
}
$ws.Range("D2").Value = "26.737.52"
$ws.Range("E2").Value = "  -2.53%  "
$ws.Range("D3").Value = "1.560.40"
$ws.Range("E3").Value = "  -0.15%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("E5").Value = "  -1.03%  "
$ws.Range("E6").Value = "  -2.08%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  +0.37%  "
$ws.Range("E9").Value = "  -0.46%  "
$ws.Range("E10").Value = "  -1.07%  "
$ws.Range("E11").Value = "  -0.21%  "
$ws.Range("D12").Value = "1.783.33"
$ws.Range("E12").Value = "  -0.11%  "
$ws.Range("D13").Value = "1.571.69"
$ws.Range("E13").Value = "  +0.56%  "
$ws.Range("D14").Value = "3.73"
$ws.Range("E14").Value = "  -2.14%  "
$ws.Range("E15").Value = "  -1.00%  "
$ws.Range("D16").Value = "61.49"
$ws.Range("E16").Value = "  -2.85%  "
$ws.Range("D17").Value = "26.757.38"
$ws.Range("E17").Value = "  -2.41%  "
$ws.Range("D18").Value = "214.19"
$ws.Range("E18").Value = "  +0.56%  "
$ws.Range("D19").Value = "7.35"
$ws.Range("E19").Value = "  +1.48%  "
$ws.Range("D20").Value = "0.0₃0676"
$ws.Range("E20").Value = "  -1.61%  "
$ws.Range("E21").Value = "  +0.08%  "
$ws.Range("E22").Value = "  -0.61%  "
$ws.Range("D23").Value = "9.33"
$ws.Range("E23").Value = "  -2.23%  "
$ws.Range("E24").Value = "  -0.47%  "
$ws.Range("D25").Value = "152.69"
$ws.Range("E25").Value = "  -0.11%  "
$ws.Range("D26").Value = "6.76"
$ws.Range("E26").Value = "  +0.46%  "
$ws.Range("D27").Value = "14.85"
$ws.Range("E27").Value = "  -0.91%  "
$ws.Range("E28").Value = "  +0.07%  "
$ws.Range("E29").Value = "  -1.21%  "
$ws.Range("E30").Value = "  -3.61%  "
$ws.Range("E31").Value = "  -1.51%  "
$ws.Range("E32").Value = "  -1.14%  "
$ws.Range("D33").Value = "1.388.35"
$ws.Range("E33").Value = "  +2.01%  "
$ws.Range("D34").Value = "2.90"
$ws.Range("E34").Value = "  -1.41%  "
$ws.Range("D35").Value = "1.56"
$ws.Range("E35").Value = "  +1.94%  "
$ws.Range("E36").Value = "  -0.57%  "
$ws.Range("D37").Value = "0.928"
$ws.Range("E37").Value = "  -4.71%  "
$ws.Range("E38").Value = "  -2.81%  "
$ws.Range("E39").Value = "  -1.51%  "
$ws.Range("E40").Value = "  -3.32%  "
$ws.Range("E41").Value = "  +0.09%  "
$ws.Range("D42").Value = "0.998"
$ws.Range("E42").Value = "  +2.39%  "
$ws.Range("D43").Value = "5.41"
$ws.Range("E43").Value = "  +2.62%  "
$ws.Range("E44").Value = "  +1.87%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").Value = "1.76"
$ws.Range("E45").Value = "  -1.28%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "63.23"
$ws.Range("E46").Value = "  -1.23%  "
$ws.Range("D47").Value = "1.695.19"
$ws.Range("E47").Value = "  -0.28%  "
$ws.Range("D48").Value = "85.66"
$ws.Range("E48").Value = "  +0.32%  "
$ws.Range("E49").Value = "  -0.05%  "
$ws.Range("E50").Value = "  -0.96%  "
$ws.Range("E51").Value = "  -0.09%  "
